$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New exercise columns -------------------------------------------------
# Insert the three new shared strings in the same order they first appear
# in the target workbook (E1 "v-crunch", F1 "penguins", D1 "shrugs") so the
# shared-strings table is built up in the same order as the authored file.
$ws.Range("E1").Value = "v-crunch"
$ws.Range("F1").Value = "penguins"
$ws.Range("D1").Value = "shrugs"

# Fill the new columns with 0 for every existing day of data (rows 2-25).
$ws.Range("D2:F25").Value = 0

# --- New day of data (row 26) ---------------------------------------------
# Copy the date formatting from the row above so the new date cell reuses
# the existing date style instead of creating a new one, then set its
# serial value explicitly (2021-02-13).
$ws.Range("A25").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A26").Value = 44240
$ws.Range("B26").Value = 30
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 40
$ws.Range("E26").Value = 35
$ws.Range("F26").Value = 150

# --- View / selection tweaks ------------------------------------------------
$ws.Range("I17").Select()
